# LeetCode workbook update:
# - Added two new problems (Transpose Matrix, Remove Duplicates from Sorted List)
# - Renamed "Running Sum of 1d Array" -> "Running Sum of 1d Array (June Day1)"
# - Added category labels for rows that were missing them (13. Linked List, 14. Arrays)
# - Added hyperlinks + hyperlink styling for rows 11-18 in column D
# - Updated the active selection cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing category (column A) labels -----------------------------
$ws.Range("A13").Value = "13. Linked List"
$ws.Range("A14").Value = "14. Arrays"

# --- Row 16: rename the existing "Running Sum of 1d Array" entry ------------
$ws.Range("B16").Value = "Running Sum of 1d Array (June Day1)"

# --- Row 17 (new): Transpose Matrix ------------------------------------------
$ws.Range("A17").Value = "17.Arrays"
$ws.Range("B17").Value = "Transpose Matrix (June Day2)"
$ws.Range("B17").Style = "Good"

# --- Row 18 (new): Remove Duplicates from Sorted List ------------------------
$ws.Range("A18").Value = "18. Linked List"
$ws.Range("B18").Value = "Remove Duplicates from Sorted List"
$ws.Range("B18").Style = "Good"

# --- Hyperlinks for column D (rows 11-18), added bottom-up to match rId order
$ws.Hyperlinks.Add($ws.Range("D18"), "https://leetcode.com/problems/remove-duplicates-from-sorted-list/")
$ws.Hyperlinks.Add($ws.Range("D17"), "https://leetcode.com/problems/transpose-matrix/")
$ws.Hyperlinks.Add($ws.Range("D16"), "https://leetcode.com/problems/running-sum-of-1d-array/")
$ws.Hyperlinks.Add($ws.Range("D15"), "https://leetcode.com/problems/longest-palindromic-substring/")
$ws.Hyperlinks.Add($ws.Range("D14"), "https://leetcode.com/problems/maximum-subarray/")
$ws.Hyperlinks.Add($ws.Range("D13"), "https://leetcode.com/problems/merge-two-sorted-lists/")
$ws.Hyperlinks.Add($ws.Range("D12"), "https://leetcode.com/problems/roman-to-integer/")
$ws.Hyperlinks.Add($ws.Range("D11"), "https://leetcode.com/problems/palindrome-number/")

# Re-apply the workbook's existing "Hyperlink" cell style (keeps the same
# style index used by D3:D10 instead of leaving the auto-generated one from
# Hyperlinks.Add on these cells).
$ws.Range("D11:D18").Style = "Hyperlink"

# --- Selection ----------------------------------------------------------------
$ws.Range("C19").Select()
